$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the Neo4j query text in cell B2 (ParticipantsTab row):
#  - drop the coalesce(... "Not specified in data") wrapper around samp.sample_id
#  - lowercase "limit" in the final ORDER BY ... LIMIT 100 clause
$oldQuery = $ws.Range("B2").Value()
$newQuery = $oldQuery.Replace(
    "WITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, ""Not specified in data""))) as samp",
    "WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp"
)
$newQuery = $newQuery.Replace(
    "ORDER BY p.participant_id LIMIT 100",
    "ORDER BY p.participant_id limit 100"
)
$ws.Range("B2").Value = $newQuery

# Move the active cell selection from D3 to D2
$ws.Range("D2").Select()
